$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "30.365.70"
$ws.Range("E2").Value = "  -0.96%  "

Set-TextValue $ws.Range("D3") "1.859.50"
$ws.Range("E3").Value = "  -1.03%  "

Set-TextValue $ws.Range("D4") "0.9998"
$ws.Range("E4").Value = "  -0.19%  "

Set-TextValue $ws.Range("D5") "234.82"
$ws.Range("E5").Value = "  -1.69%  "

$ws.Range("E6").Value = "  -0.14%  "

Set-TextValue $ws.Range("D7") "0.4740"
$ws.Range("E7").Value = "  -1.47%  "

Set-TextValue $ws.Range("D8") "0.2746"
$ws.Range("E8").Value = "  -2.97%  "

Set-TextValue $ws.Range("D9") "0.06443"
$ws.Range("E9").Value = "  -1.17%  "

Set-TextValue $ws.Range("D10") "1.875.79"
$ws.Range("E10").Value = "  -1.80%  "

Set-TextValue $ws.Range("D11") "0.07435"
$ws.Range("E11").Value = "  -0.66%  "

Set-TextValue $ws.Range("D12") "16.35"
$ws.Range("E12").Value = "  -1.12%  "

Set-TextValue $ws.Range("D13") "5.010"
$ws.Range("E13").Value = "  -1.65%  "

Set-TextValue $ws.Range("D14") "85.70"
$ws.Range("E14").Value = "  -2.78%  "

Set-TextValue $ws.Range("D15") "0.6363"
$ws.Range("E15").Value = "  -4.18%  "

Set-TextValue $ws.Range("D16") "30.331.42"
$ws.Range("E16").Value = "  -0.92%  "

Set-TextValue $ws.Range("D17") "0.9997"
$ws.Range("E17").Value = "  -0.13%  "

Set-TextValue $ws.Range("D18") "231.75"
$ws.Range("E18").Value = "  +1.36%  "

Set-TextValue $ws.Range("D19") "12.80"
$ws.Range("E19").Value = "  -3.61%  "

Set-TextValue $ws.Range("D20") "0.000007429"
$ws.Range("E20").Value = "  -2.16%  "

Set-TextValue $ws.Range("D21") "2.099.91"
$ws.Range("E21").Value = "  -4.17%  "

Set-TextValue $ws.Range("D22") "1.001"
$ws.Range("E22").Value = "  -0.11%  "

Set-TextValue $ws.Range("D23") "5.018"
$ws.Range("E23").Value = "  -4.96%  "

Set-TextValue $ws.Range("D24") "6.016"
$ws.Range("E24").Value = "  -2.24%  "

Set-TextValue $ws.Range("D25") "9.299"
$ws.Range("E25").Value = "  +0.21%  "

Set-TextValue $ws.Range("D26") "165.61"
$ws.Range("E26").Value = "  -1.90%  "

$ws.Range("E27").Value = "  -2.99%  "

$ws.Range("E28").Value = "  -1.81%  "

Set-TextValue $ws.Range("D29") "0.1042"
$ws.Range("E29").Value = "  +7.39%  "

Set-TextValue $ws.Range("D30") "1.391"
$ws.Range("E30").Value = "  -0.96%  "

Set-TextValue $ws.Range("D31") "4.149"
$ws.Range("E31").Value = "  -4.42%  "

Set-TextValue $ws.Range("D32") "3.938"
$ws.Range("E32").Value = "  -1.72%  "

Set-TextValue $ws.Range("D33") "0.04911"
$ws.Range("E33").Value = "  -3.30%  "

Set-TextValue $ws.Range("D34") "1.153"
$ws.Range("E34").Value = "  -5.48%  "

Set-TextValue $ws.Range("D35") "0.7277"
$ws.Range("E35").Value = "  -2.99%  "

Set-TextValue $ws.Range("D36") "0.9993"
$ws.Range("E36").Value = "  -0.58%  "

Set-TextValue $ws.Range("D37") "2.696"
$ws.Range("E37").Value = "  -0.65%  "

Set-TextValue $ws.Range("D38") "0.01910"
$ws.Range("E38").Value = "  +2.52%  "

Set-TextValue $ws.Range("D39") "2.651"
$ws.Range("E39").Value = "  +0.42%  "

Set-TextValue $ws.Range("D40") "0.9111"
$ws.Range("E40").Value = "  -0.24%  "

Set-TextValue $ws.Range("D41") "1.974"
$ws.Range("E41").Value = "  -5.09%  "

Set-TextValue $ws.Range("D42") "105.59"
$ws.Range("E42").Value = "  -0.63%  "

Set-TextValue $ws.Range("D43") "0.9996"
$ws.Range("E43").Value = "  -0.01%  "

Set-TextValue $ws.Range("D44") "0.4124"
$ws.Range("E44").Value = "  -3.45%  "

Set-TextValue $ws.Range("D45") "5.567"
$ws.Range("E45").Value = "  -3.47%  "

Set-TextValue $ws.Range("D46") "7.157"
$ws.Range("E46").Value = "  -2.46%  "

Set-TextValue $ws.Range("D47") "61.29"
$ws.Range("E47").Value = "  -4.66%  "

Set-TextValue $ws.Range("D48") "0.1214"
$ws.Range("E48").Value = "  -5.90%  "

Set-TextValue $ws.Range("D49") "8.737"
$ws.Range("E49").Value = "  -2.76%  "

Set-TextValue $ws.Range("D50") "1.410"
$ws.Range("E50").Value = "  -4.41%  "

Set-TextValue $ws.Range("D51") "33.44"
$ws.Range("E51").Value = "  -0.86%  "

